# Update forest data - 2026-01-07 12:21
#
# The "New" sheet holds 4 freshly-scraped listings (rows 2-5). This pass
# files them into the bottom of "Previously added" (rows 360-363) and
# replaces "New" with a single newly-scraped listing.

$wb = $excel.ActiveWorkbook
$wsPrev = $wb.Worksheets.Item("Previously added")
$wsNew  = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# 1) Capture the 4 rows currently sitting in "New" (rows 2..5) along with
#    the hyperlink target on column A of each, before anything moves.
# ---------------------------------------------------------------------
$moveRows = @()
for ($r = 2; $r -le 5; $r++) {
    $link = $null
    for ($i = 1; $i -le $wsNew.Hyperlinks.Count; $i++) {
        $hl = $wsNew.Hyperlinks.Item($i)
        if ($hl.Range.Row -eq $r) { $link = $hl.Address }
    }
    $moveRows += , @(
        $wsNew.Cells.Item($r, 1).Value2,
        $wsNew.Cells.Item($r, 2).Value2,
        $wsNew.Cells.Item($r, 3).Value2,
        $wsNew.Cells.Item($r, 4).Value2,
        $wsNew.Cells.Item($r, 5).Value2,
        $wsNew.Cells.Item($r, 6).Value2,
        $link
    )
}

# ---------------------------------------------------------------------
# 2) Append those 4 rows to the bottom of "Previously added".
# ---------------------------------------------------------------------
$destFirst = $wsPrev.UsedRange.Rows.Count + 1

for ($k = 0; $k -lt $moveRows.Count; $k++) {
    $destRow = $destFirst + $k
    $data = $moveRows[$k]

    $wsPrev.Cells.Item($destRow, 1).Value = $data[0]
    $wsPrev.Cells.Item($destRow, 2).Value = $data[1]
    $wsPrev.Cells.Item($destRow, 3).Value = $data[2]
    $wsPrev.Cells.Item($destRow, 4).Value = $data[3]
    $wsPrev.Cells.Item($destRow, 5).Value = $data[4]
    $wsPrev.Cells.Item($destRow, 6).Value = $data[5]

    if ($data[6]) {
        $wsPrev.Hyperlinks.Add($wsPrev.Cells.Item($destRow, 1), $data[6])
        # Hyperlinks.Add stamps the built-in "Hyperlink" style; restore the
        # sheet's normal link-column look by copying formats from the row
        # right above (already styled like every other link cell).
        $wsPrev.Cells.Item($destRow - 1, 1).Copy()
        $wsPrev.Cells.Item($destRow, 1).PasteSpecial(-4122)  # xlPasteFormats
    }
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3) Drop the now-migrated rows 3..5 from "New", keeping row 2 as the
#    template for the newly scraped listing, then blank out every
#    left-behind hyperlink on the sheet (row deletion does not clean
#    those up by itself).
# ---------------------------------------------------------------------
$wsNew.Rows.Item(5).Delete()
$wsNew.Rows.Item(4).Delete()
$wsNew.Rows.Item(3).Delete()
$wsNew.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 4) Write the single new listing into row 2 of "New".
# ---------------------------------------------------------------------
$newLink = "https://www.ss.com/msg/lv/real-estate/wood/other/deobg.html"

$wsNew.Cells.Item(2, 1).Value = $newLink
$wsNew.Cells.Item(2, 2).Value = "9 000 €"
$wsNew.Cells.Item(2, 3).Value = ""
$wsNew.Cells.Item(2, 4).Value = "4 ha."
$wsNew.Cells.Item(2, 5).Value = "42760030110"
$wsNew.Cells.Item(2, 6).Value = 46028.9125

$wsNew.Hyperlinks.Add($wsNew.Range("A2"), $newLink)

# Restore A2's original (non built-in) hyperlink-column style, the same
# way we fixed up the moved rows above.
$wsPrev.Cells.Item(2, 1).Copy()
$wsNew.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
